$d = $word.ActiveDocument

# 1) Remove the stray _GoBack bookmark that currently sits in the
#    "Link do chat" paragraph (right after the hyperlink run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Make the "Zero Data: iniciando sua jornada no mundo dos dados;"
#    bullet bold (paragraph mark + both runs).
$target = $d.Paragraphs.Item(11)
$target.Range.Font.Bold = 1

# 3) Re-add the _GoBack bookmark, now wrapping just the title run
#    ("Zero Data: iniciando sua jornada no mundo dos dados"), not the
#    trailing semicolon run.
$titleText = "Zero Data: iniciando sua jornada no mundo dos dados"
$start = $target.Range.Start
$bookmarkRange = $d.Range($start, $start + $titleText.Length)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
